$d = $word.ActiveDocument

$pairs = @(
    @("975÷2=487, 1", "804÷8=100, 4"),
    @("779÷7=111, 2", "388÷7=55, 3"),
    @("662÷7=94, 4", "196÷8=24, 4"),
    @("909÷2=454, 1", "908÷6=151, 2"),
    @("117÷7=16, 5", "162÷8=20, 2"),
    @("176÷3=58, 2", "854÷3=284, 2"),
    @("454÷2=227, 0", "277÷4=69, 1"),
    @("966÷9=107, 3", "387÷3=129, 0"),
    @("859÷6=143, 1", "878÷7=125, 3"),
    @("166÷9=18, 4", "444÷3=148, 0"),
    @("534÷6=89, 0", "503÷6=83, 5"),
    @("606÷6=101, 0", "127÷3=42, 1"),
    @("534÷7=76, 2", "923÷8=115, 3"),
    @("791÷9=87, 8", "120÷8=15, 0"),
    @("776÷9=86, 2", "842÷6=140, 2"),
    @("885÷5=177, 0", "767÷2=383, 1"),
    @("346÷3=115, 1", "710÷3=236, 2"),
    @("651÷2=325, 1", "539÷9=59, 8"),
    @("343÷3=114, 1", "477÷6=79, 3"),
    @("259÷3=86, 1", "430÷3=143, 1"),
    @("112÷6=18, 4", "839÷4=209, 3"),
    @("411÷3=137, 0", "892÷8=111, 4"),
    @("719÷8=89, 7", "343÷9=38, 1"),
    @("266÷8=33, 2", "398÷7=56, 6"),
    @("533÷5=106, 3", "519÷7=74, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
